$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 15: Objetos | Item,Categorizable | La clase Item hereda de la clase categorizable. | (blank) | Creamos la clase "categorizable" ... ---
$ws.Range("A15").Value = "Objetos"
$ws.Range("B15").Value = "Item,Categorizable"
$ws.Range("C15").Value = "La clase Item hereda de la clase categorizable."
$ws.Range("E15").Value = "Creamos la clase ""categorizable"" para que todos los objetos que tengan que categorizarse utilize los atributos y métodos de esta clase, delegando la responsabilidad de de asignar categorias a los items."

# --- Row 16: Objetos | Entidad | Colocamos atributo criterios ... | (blank) | (blank) ---
$ws.Range("A16").Value = "Objetos"
$ws.Range("B16").Value = "Entidad"
$ws.Range("C16").Value = "Colocamos atributo criterios para llevar registro de los criterios adoptados por la organización"

# --- Row 17: Objetos | Egreso,Ingreso | Relacionamos ambos ... ---
$ws.Range("A17").Value = "Objetos"
$ws.Range("B17").Value = "Egreso,Ingreso"
$ws.Range("C17").Value = "Relacionamos ambos colocando a cada uno como atributo del otro"

# --- Row 18: Objetos | Categoria | La clase categoria conoce al criterio ... | (blank) | De esta forma cada Item ... ---
$ws.Range("A18").Value = "Objetos"
$ws.Range("B18").Value = "Categoria"
$ws.Range("C18").Value = "La clase categoria conoce al criterio que pertenece"
$ws.Range("E18").Value = "De esta forma cada Item (que hereda de Categorizable) tiene una lista de categorias y a su vez tambien conoce al criterio que pertenece."

# --- Row 19: Objetos | Criterio, Entidad | El orden de la lista representa la jerarquia | (blank) | El usuario ingresa los criterios ... ---
$ws.Range("A19").Value = "Objetos"
$ws.Range("B19").Value = "Criterio, Entidad"
$ws.Range("C19").Value = "El orden de la lista representa la jerarquia"
$ws.Range("E19").Value = "El usuario ingresa los criterios en el orden de jerarquia"

# --- Row 20: Objetos | Criterio, Categorizable | Ordenar la lista de categorias ... | (blank) | De esta forma cada item conoce su jerarquia ... ---
$ws.Range("A20").Value = "Objetos"
$ws.Range("B20").Value = "Criterio, Categorizable"
$ws.Range("C20").Value = "Ordenar la lista de categorias según la jerarquia de los criterios"
$ws.Range("E20").Value = "De esta forma cada item conoce su jerarquia, porque conoce su criterio. "

# --- Row 17 justificacion (edited last, fixing "perteneco" -> "pertenece" typo) ---
$ws.Range("E17").Value = "Para asociar a ambos y que sea mas facil llevar registro de que ingreso pertenece a que egreso y viceversa"

# --- Update the selected/visible cell to reflect the new view state ---
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("E17").Select()
